# Auto-generated by inspection of the OOXML diff for Faerie_Profits workbook.
# Updates cached market-data values (columns H..N) on several rows across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1004
$ws.Range("I10").Value = 1004
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1004
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -711
$ws.Range("N10").Value = ""

$ws.Range("H80").Value = 1609.4
$ws.Range("I80").Value = 979.85
$ws.Range("J80").Value = 2113.04
$ws.Range("K80").Value = 2939.55
$ws.Range("L80").Value = 6339.12
$ws.Range("M80").Value = -1941.55
$ws.Range("N80").Value = -8335.119999999999

$ws.Range("H83").Value = 1609.4
$ws.Range("I83").Value = 979.85
$ws.Range("J83").Value = 2113.04
$ws.Range("K83").Value = 8818.65
$ws.Range("L83").Value = 19017.36
$ws.Range("M83").Value = -3826.65
$ws.Range("N83").Value = -29001.36

$ws.Range("H138").Value = 234384.95
$ws.Range("J138").Value = 418962.6
$ws.Range("L138").Value = 1256887.8
$ws.Range("N138").Value = -1267167.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 5956.778
$ws.Range("I5").Value = 5956.778
$ws.Range("K5").Value = 5956.778
$ws.Range("M5").Value = -5844.778

$ws.Range("H14").Value = 4756.25
$ws.Range("I14").Value = 749.5
$ws.Range("K14").Value = 749.5
$ws.Range("M14").Value = -574.5

$ws.Range("H45").Value = 4090.0833
$ws.Range("I45").Value = 2010.6364
$ws.Range("K45").Value = 2010.6364
$ws.Range("M45").Value = -1633.6364

$ws.Range("H62").Value = 575059.75
$ws.Range("J62").Value = 575059.75
$ws.Range("L62").Value = 575059.75
$ws.Range("N62").Value = -576307.75

$ws.Range("H65").Value = 575059.75
$ws.Range("J65").Value = 575059.75
$ws.Range("L65").Value = 1725179.25
$ws.Range("N65").Value = -1731419.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 5956.778
$ws.Range("I4").Value = 5956.778
$ws.Range("K4").Value = 5956.778
$ws.Range("M4").Value = -5841.778

$ws.Range("H80").Value = 560.13336
$ws.Range("I80").Value = 290.4
$ws.Range("J80").Value = 695
$ws.Range("K80").Value = 290.4
$ws.Range("L80").Value = 695
$ws.Range("M80").Value = 707.6
$ws.Range("N80").Value = -2691

$ws.Range("H83").Value = 560.13336
$ws.Range("I83").Value = 290.4
$ws.Range("J83").Value = 695
$ws.Range("K83").Value = 1452
$ws.Range("L83").Value = 3475
$ws.Range("M83").Value = 3540
$ws.Range("N83").Value = -13459

$ws.Range("H86").Value = 3333
$ws.Range("I86").Value = 3333
$ws.Range("K86").Value = 3333
$ws.Range("M86").Value = -2210

$ws.Range("H89").Value = 3333
$ws.Range("I89").Value = 3333
$ws.Range("K89").Value = 16665
$ws.Range("M89").Value = -11049

$ws.Range("H99").Value = 2322.111
$ws.Range("I99").Value = 1984.1428
$ws.Range("K99").Value = 1984.1428
$ws.Range("M99").Value = -486.1428000000001

$ws.Range("H134").Value = 2946.5522
$ws.Range("I134").Value = 2602.5417
$ws.Range("K134").Value = 7807.625100000001
$ws.Range("M134").Value = -5272.625100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 836.4
$ws.Range("I22").Value = 848.75
$ws.Range("J22").Value = 787
$ws.Range("K22").Value = 848.75
$ws.Range("L22").Value = 787
$ws.Range("M22").Value = -498.75
$ws.Range("N22").Value = -1487

$ws.Range("H31").Value = 1937.9546
$ws.Range("I31").Value = 1402
$ws.Range("K31").Value = 1402
$ws.Range("M31").Value = -1107

$ws.Range("H34").Value = 1937.9546
$ws.Range("I34").Value = 1402
$ws.Range("K34").Value = 1402
$ws.Range("M34").Value = -1200

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 131.77777
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 131.77777
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 790.66662
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = -1016.66662

$ws.Range("H50").Value = 1109
$ws.Range("J50").Value = 1000
$ws.Range("L50").Value = 3000
$ws.Range("N50").Value = -3962

$ws.Range("H53").Value = 1109
$ws.Range("J53").Value = 1000
$ws.Range("L53").Value = 3000
$ws.Range("N53").Value = -3962

$ws.Range("H80").Value = 3798.4
$ws.Range("J80").Value = 5664
$ws.Range("L80").Value = 16992
$ws.Range("N80").Value = -18864

$ws.Range("H83").Value = 3798.4
$ws.Range("J83").Value = 5664
$ws.Range("L83").Value = 50976
$ws.Range("N83").Value = -60336

$ws.Range("H129").Value = 1601.7778
$ws.Range("I129").Value = 1033.5834
$ws.Range("J129").Value = 2738.1667
$ws.Range("K129").Value = 3100.7502
$ws.Range("L129").Value = 8214.500100000001
$ws.Range("M129").Value = 1899.2498
$ws.Range("N129").Value = -18214.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 31502.666
$ws.Range("J20").Value = 31502.666
$ws.Range("L20").Value = 31502.666
$ws.Range("N20").Value = -31992.666

$ws.Range("H21").Value = 17776.5
$ws.Range("I21").Value = 9998
$ws.Range("J21").Value = 25555
$ws.Range("K21").Value = 9998
$ws.Range("L21").Value = 25555
$ws.Range("M21").Value = -9825
$ws.Range("N21").Value = -25901

$ws.Range("H24").Value = 22508752
$ws.Range("I24").Value = 45005504
$ws.Range("K24").Value = 45005504
$ws.Range("M24").Value = -45005331

$ws.Range("H30").Value = 17776.5
$ws.Range("I30").Value = 9998
$ws.Range("J30").Value = 25555
$ws.Range("K30").Value = 9998
$ws.Range("L30").Value = 25555
$ws.Range("M30").Value = -9893
$ws.Range("N30").Value = -25765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3172.5334
$ws.Range("J22").Value = 3773
$ws.Range("L22").Value = 3773
$ws.Range("N22").Value = -4363

$ws.Range("H27").Value = 3172.5334
$ws.Range("J27").Value = 3773
$ws.Range("L27").Value = 3773
$ws.Range("N27").Value = -3987

$ws.Range("H33").Value = 18853.5
$ws.Range("I33").Value = 18853.5
$ws.Range("K33").Value = 18853.5
$ws.Range("M33").Value = -18563.5

$ws.Range("H40").Value = 7271.385
$ws.Range("I40").Value = 7312.75
$ws.Range("K40").Value = 7312.75
$ws.Range("M40").Value = -7176.75

$ws.Range("H61").Value = 3599.8
$ws.Range("I61").Value = 3499.75
$ws.Range("K61").Value = 3499.75
$ws.Range("M61").Value = -3297.75

$ws.Range("H62").Value = 89991.664
$ws.Range("J62").Value = 89991.664
$ws.Range("L62").Value = 89991.664
$ws.Range("N62").Value = -91239.664

$ws.Range("H64").Value = 74996.664
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 74996.664
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 74996.664
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -75446.664

$ws.Range("H65").Value = 89991.664
$ws.Range("J65").Value = 89991.664
$ws.Range("L65").Value = 269974.992
$ws.Range("N65").Value = -276214.992

$ws.Range("H67").Value = 74996.664
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 74996.664
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 74996.664
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -76556.664

$ws.Range("H68").Value = 5832.6665
$ws.Range("I68").Value = 6250.25
$ws.Range("K68").Value = 6250.25
$ws.Range("M68").Value = -5501.25

$ws.Range("H71").Value = 5832.6665
$ws.Range("I71").Value = 6250.25
$ws.Range("K71").Value = 31251.25
$ws.Range("M71").Value = -27507.25

$ws.Range("H113").Value = 3599.8
$ws.Range("I113").Value = 3499.75
$ws.Range("K113").Value = 3499.75
$ws.Range("M113").Value = -1329.75

$ws.Range("H122").Value = 13061.8125
$ws.Range("I122").Value = 9153.308000000001
$ws.Range("K122").Value = 27459.924
$ws.Range("M122").Value = -25009.924

$ws.Range("H132").Value = 1756.5
$ws.Range("J132").Value = 499
$ws.Range("L132").Value = 1497
$ws.Range("N132").Value = -6557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 976
$ws.Range("J22").Value = 1015
$ws.Range("L22").Value = 1015
$ws.Range("N22").Value = -1601

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""

$ws.Range("H29").Value = 2375
$ws.Range("J29").Value = 2500
$ws.Range("L29").Value = 2500
$ws.Range("N29").Value = -3080

$ws.Range("H37").Value = 44079.668
$ws.Range("I37").Value = 27419.5
$ws.Range("J37").Value = 77400
$ws.Range("K37").Value = 27419.5
$ws.Range("L37").Value = 77400
$ws.Range("M37").Value = -27216.5
$ws.Range("N37").Value = -77806

$ws.Range("H113").Value = 8334916
$ws.Range("I113").Value = 83333336
$ws.Range("J113").Value = 1758.4445
$ws.Range("K113").Value = 250000008
$ws.Range("L113").Value = 5275.333500000001
$ws.Range("M113").Value = -249997838
$ws.Range("N113").Value = -9615.333500000001
